$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue ($ws.Range("D2")) "251.62"
Set-TextValue ($ws.Range("D3")) "23.00"
Set-TextValue ($ws.Range("D4")) "5.447"
Set-TextValue ($ws.Range("D5")) "0.05664"
Set-TextValue ($ws.Range("D6")) "3.441"
Set-TextValue ($ws.Range("D7")) "6.389"
Set-TextValue ($ws.Range("D8")) "0.8198"
Set-TextValue ($ws.Range("D9")) "0.9303"
Set-TextValue ($ws.Range("D10")) "0.1438"
Set-TextValue ($ws.Range("D11")) "0.07489"
Set-TextValue ($ws.Range("D12")) "0.03164"
Set-TextValue ($ws.Range("D13")) "0.03076"
Set-TextValue ($ws.Range("D14")) "0.09353"
Set-TextValue ($ws.Range("D15")) "3.556"
Set-TextValue ($ws.Range("D16")) "0.001641"
Set-TextValue ($ws.Range("D17")) "0.04726"
Set-TextValue ($ws.Range("D18")) "0.0005777"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue ($ws.Range("D19")) "0.006347"
Set-TextValue ($ws.Range("D20")) "0.005037"
Set-TextValue ($ws.Range("D21")) "0.001030"
Set-TextValue ($ws.Range("D22")) "0.0001497"
Set-TextValue ($ws.Range("D24")) "2.198"
Set-TextValue ($ws.Range("D25")) "0.3293"
$ws.Range("E27").Value = "26AAXTokenAAB"
Set-TextValue ($ws.Range("D28")) "0.0002993"
Set-TextValue ($ws.Range("D40")) "0.04018"
Set-TextValue ($ws.Range("D41")) "0.006936"
Set-TextValue ($ws.Range("D42")) "0.1074"
Set-TextValue ($ws.Range("D43")) "0.002778"
Set-TextValue ($ws.Range("D44")) "0.007604"
Set-TextValue ($ws.Range("D45")) "0.00005560"
Set-TextValue ($ws.Range("D46")) "0.00000000749"
Set-TextValue ($ws.Range("D48")) "0.6585"
Set-TextValue ($ws.Range("D49")) "0.2239"
Set-TextValue ($ws.Range("D50")) "0.00002096"
Set-TextValue ($ws.Range("D51")) "0.01008"
